# Refresh the "Price" (column D) and "Volume(1h)" (column E) columns of the
# cryptos sheet with the latest scraped snapshot.
#
# Column D values are plain text (e.g. "55.386.96", "0.999") rather than
# numbers, so a leading apostrophe is used to stop Excel's input parser from
# re-interpreting digit-looking strings (like "1.00" or "0.0983") as
# numbers; the style is then reset to "Normal" so the apostrophe's
# quote-prefix flag doesn't linger as a formatting change.
# Column E values already carry padding spaces ("  +0.79%  "), which keeps
# Excel from parsing them as percentages, so they can be assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'54.931.54"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.79%  "

$ws.Cells.Item(3, 4).Value = "'2.292.31"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +0.10%  "

$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.14%  "

$ws.Cells.Item(5, 4).Value = "'505.95"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.40%  "

$ws.Cells.Item(6, 4).Value = "'129.66"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.66%  "

$ws.Cells.Item(7, 5).Value = "  -0.25%  "

$ws.Cells.Item(8, 4).Value = "'0.531"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.20%  "

$ws.Cells.Item(9, 4).Value = "'2.313.67"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.47%  "

$ws.Cells.Item(10, 4).Value = "'0.0983"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +2.55%  "

$ws.Cells.Item(12, 5).Value = "  +6.61%  "

$ws.Cells.Item(13, 5).Value = "  +0.45%  "

$ws.Cells.Item(14, 4).Value = "'23.81"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +3.84%  "

$ws.Cells.Item(15, 4).Value = "'2.702.18"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.05%  "

$ws.Cells.Item(16, 4).Value = "'54.982.20"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.92%  "

$ws.Cells.Item(17, 5).Value = "  +1.35%  "

$ws.Cells.Item(18, 4).Value = "'2.308.47"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.49%  "

$ws.Cells.Item(19, 4).Value = "'10.52"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +2.25%  "

$ws.Cells.Item(20, 5).Value = "  +0.30%  "

$ws.Cells.Item(21, 4).Value = "'310.95"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +2.05%  "

$ws.Cells.Item(22, 4).Value = "'6.61"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +3.82%  "

$ws.Cells.Item(23, 4).Value = "'0.998"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.28%  "

$ws.Cells.Item(24, 4).Value = "'60.31"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -2.67%  "

$ws.Cells.Item(25, 4).Value = "'0.995"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.22%  "

$ws.Cells.Item(26, 5).Value = "  +0.21%  "

$ws.Cells.Item(27, 4).Value = "'7.52"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +2.33%  "

$ws.Cells.Item(28, 4).Value = "'173.05"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +0.84%  "

$ws.Cells.Item(29, 4).Value = "'6.15"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +3.08%  "

$ws.Cells.Item(30, 4).Value = "'0.0₃0708"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +1.74%  "

$ws.Cells.Item(31, 5).Value = "  +0.23%  "

$ws.Cells.Item(32, 4).Value = "'1.14"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +4.53%  "

$ws.Cells.Item(33, 4).Value = "'0.998"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.04%  "

$ws.Cells.Item(34, 4).Value = "'18.04"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.10%  "

$ws.Cells.Item(35, 4).Value = "'0.994"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.23%  "

$ws.Cells.Item(36, 5).Value = "  -5.42%  "

$ws.Cells.Item(37, 5).Value = "  +2.40%  "

$ws.Cells.Item(38, 4).Value = "'3.89"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +3.84%  "

$ws.Cells.Item(39, 4).Value = "'36.83"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.57%  "

$ws.Cells.Item(40, 5).Value = "  +2.30%  "

$ws.Cells.Item(41, 4).Value = "'0.377"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.26%  "

$ws.Cells.Item(42, 4).Value = "'133.81"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +5.67%  "

$ws.Cells.Item(43, 5).Value = "  +1.10%  "

$ws.Cells.Item(44, 4).Value = "'4.93"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -1.65%  "

$ws.Cells.Item(45, 4).Value = "'259.98"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +7.09%  "

$ws.Cells.Item(46, 4).Value = "'0.0506"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.87%  "

$ws.Cells.Item(47, 5).Value = "  +1.66%  "

$ws.Cells.Item(48, 4).Value = "'0.550"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.19%  "

$ws.Cells.Item(49, 4).Value = "'0.376"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.33%  "

$ws.Cells.Item(50, 5).Value = "  +1.92%  "

$ws.Cells.Item(51, 4).Value = "'16.47"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.33%  "

